# Update Betfair Back/Lay odds figures for 2025-12-04 on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AJ2").Value = 50
$ws.Range("AL2").Value = 42
$ws.Range("AN2").Value = 24
$ws.Range("G2").Value = 2.98
$ws.Range("K2").Value = 3.7
$ws.Range("Q2").Value = 1.82
$ws.Range("V2").Value = 1.59
$ws.Range("W2").Value = 1.5
$ws.Range("H3").Value = 1.81
$ws.Range("K3").Value = 9.800000000000001
$ws.Range("Q3").Value = 1.48
$ws.Range("L4").Value = 1.31
$ws.Range("AK5").Value = 25
$ws.Range("J5").Value = 3.3
$ws.Range("K5").Value = 3.8
$ws.Range("U5").Value = 1.67
$ws.Range("Q6").Value = 1.84
$ws.Range("V6").Value = 2.92
$ws.Range("G7").Value = 3.1
$ws.Range("H7").Value = 2.8
$ws.Range("N7").Value = 2.88
$ws.Range("Q7").Value = 2.34
$ws.Range("W7").Value = 1.47
$ws.Range("K8").Value = 4.8
$ws.Range("S8").Value = 2.72
$ws.Range("F9").Value = 1.2
$ws.Range("I9").Value = 19.5
$ws.Range("K9").Value = 9.4
$ws.Range("N9").Value = 3.55
$ws.Range("P9").Value = 3.55
$ws.Range("Q9").Value = 1.26
$ws.Range("S9").Value = 1.48
$ws.Range("W9").Value = 4.3
$ws.Range("G10").Value = 2.02
$ws.Range("I10").Value = 5.3
$ws.Range("J10").Value = 3.15
$ws.Range("U10").Value = 1.74
$ws.Range("W10").Value = 1.98
$ws.Range("AC12").Value = 13
$ws.Range("AH12").Value = 17
$ws.Range("AK12").Value = 1000
$ws.Range("AN12").Value = 16.5
$ws.Range("K12").Value = 4.5
$ws.Range("L12").Value = 1.24
$ws.Range("X12").Value = 34
$ws.Range("Y12").Value = 22
$ws.Range("AE13").Value = 80
$ws.Range("AJ13").Value = 14
$ws.Range("AK13").Value = 13.5
$ws.Range("AN13").Value = 5.6
$ws.Range("J13").Value = 5.1
$ws.Range("T13").Value = 1.74
$ws.Range("U13").Value = 2.28
$ws.Range("AF14").Value = 40
$ws.Range("AG14").Value = 20
$ws.Range("AI14").Value = 27
$ws.Range("AL14").Value = 40
$ws.Range("AM14").Value = 55
$ws.Range("AN14").Value = 26
$ws.Range("AO14").Value = 10
$ws.Range("F14").Value = 3.45
$ws.Range("G14").Value = 4.2
$ws.Range("H14").Value = 1.87
$ws.Range("I14").Value = 2.08
$ws.Range("J14").Value = 4
$ws.Range("K14").Value = 4.6
$ws.Range("N14").Value = 5.8
$ws.Range("S14").Value = 2.12
$ws.Range("U14").Value = 2.56
$ws.Range("V14").Value = 1.92
$ws.Range("W14").Value = 1.32
$ws.Range("AD15").Value = 19
$ws.Range("AH15").Value = 23
$ws.Range("AM15").Value = 140
$ws.Range("F15").Value = 2.2
$ws.Range("G15").Value = 2.28
$ws.Range("H15").Value = 3.85
$ws.Range("J15").Value = 3.15
$ws.Range("K15").Value = 3.25
$ws.Range("N15").Value = 3.3
$ws.Range("O15").Value = 1.38
$ws.Range("P15").Value = 1.77
$ws.Range("T15").Value = 1.86
$ws.Range("U15").Value = 2.02
$ws.Range("W15").Value = 1.78
$ws.Range("F16").Value = 1.9
$ws.Range("G16").Value = 2.02
$ws.Range("I16").Value = 4.9
$ws.Range("K16").Value = 4.1
$ws.Range("W16").Value = 1.98
$ws.Range("G17").Value = 1.23
$ws.Range("H17").Value = 11.5
$ws.Range("K17").Value = 16.5
$ws.Range("N17").Value = 9.199999999999999
$ws.Range("P17").Value = 3.55
$ws.Range("N18").Value = 1.1
$ws.Range("S18").Value = 3.4
$ws.Range("U18").Value = 1.04
